# Simplify NB01 to linear teaching flow and refresh regression datasets
#
# This script refreshes the "y" column (column B) values for a subset of
# rows in the regression dataset found on the active worksheet, matching
# the updated dataset values produced after the notebook's data-refresh
# step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B39").Value = 147.5846604714236
$ws.Range("B40").Value = 150.3752924944157
$ws.Range("B41").Value = 146.9029309696913
$ws.Range("B42").Value = 137.2397543573069

$ws.Range("B55").Value = 130.5187793558437
$ws.Range("B56").Value = 118.6943887251082
$ws.Range("B57").Value = 104.9372488363266
$ws.Range("B58").Value = 110.7328996911328
$ws.Range("B59").Value = 98.99442956281302
$ws.Range("B60").Value = 113.0225478581135
$ws.Range("B61").Value = 124.6236344154317
$ws.Range("B62").Value = 145.963580174758
$ws.Range("B63").Value = 145.4050488910343

$ws.Range("B72").Value = 110.8893995203106
$ws.Range("B73").Value = 117.3294001362893
